# [Maintain_Job_Status] : Test Data Added
# Update the test data row (row 8) on the "jobStatus" sheet with new
# sample values, and move the active selection to H8 with the view
# scrolled back so column A is visible again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Job No
$ws.Range("C8").Value = "CHN/BKG/AFE/00195/23-24"
# Status Date
$ws.Range("E8").Value = "10-Jul-2023"
# Completion Date
$ws.Range("J8").Value = "14-Jul-2023"
# Remarks - Hour (HH) - keep as quoted text so the list-validated
# zero-padded value ("15") is not reinterpreted as a number.
$ws.Range("F8").Value = "'15"
# Remarks - Minute (MM)
$ws.Range("G8").Value = "'06"

# Move the active selection to H8. Selecting a cell also resets the
# sheet's scrolled top-left cell back to its default (A1), matching the
# removal of the explicit topLeftCell="B1" view setting.
$ws.Activate()
$ws.Range("H8").Select()
